$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B21").Value = 6242
$ws.Range("D21").Value = 5616943
$ws.Range("E21").Value = 899.8627042614546
$ws.Range("F21").Value = 8.34924492275646
$ws.Range("H21").Value = 28.21225517990072
